# Refresh cryptos list (prices / 1h volume %) per the Jul 15 2024 GitHub Actions run.
# D/E columns hold plain text (not numbers), so any cell whose new value parses as a
# pure number gets NumberFormat "@" (Text) first -- mirrors what real Excel COM
# automation does to keep e.g. "153.74" from being auto-coerced into a float.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '62.983.31'
$ws.Range('E2').Value = '  +4.96%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.362.58'
$ws.Range('E3').Value = '  +5.25%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.00%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '561.24'
$ws.Range('E5').Value = '  +4.28%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.74'
$ws.Range('E6').Value = '  +6.07%  '

# Row 7: USDC
$ws.Range('E7').Value = '  -0.01%  '

# Row 8: XRP
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.531'
$ws.Range('E8').Value = '  +0.79%  '

# Row 9: Toncoin
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '7.55'
$ws.Range('E9').Value = '  +2.76%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  +4.52%  '

# Row 11: Cardano
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.440'
$ws.Range('E11').Value = '  +2.10%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '3.937.63'
$ws.Range('E12').Value = '  +5.12%  '

# Row 13: TRON
$ws.Range('E13').Value = '  +0.39%  '

# Row 14: ShibaInu
$ws.Range('E14').Value = '  +4.10%  '

# Row 15: Avalanche
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '27.19'
$ws.Range('E15').Value = '  +4.41%  '

# Row 16: WrappedBTC
$ws.Range('D16').Value = '63.043.78'
$ws.Range('E16').Value = '  +5.00%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '3.352.09'
$ws.Range('E17').Value = '  +4.76%  '

# Row 18: Polkadot
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.52'
$ws.Range('E18').Value = '  +4.59%  '

# Row 19: Chainlink
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.85'
$ws.Range('E19').Value = '  +5.73%  '

# Row 20: Uniswap
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '8.48'
$ws.Range('E20').Value = '  +1.49%  '

# Row 21: BitcoinCash
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '389.75'
$ws.Range('E21').Value = '  +1.74%  '

# Row 22: Polygon
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.544'
$ws.Range('E22').Value = '  +2.59%  '

# Row 23: Dai
$ws.Range('E23').Value = '  +0.21%  '

# Row 24: Litecoin
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.59'
$ws.Range('E24').Value = '  +0.41%  '

# Row 25: Kaspa
$ws.Range('E25').Value = '  +5.15%  '

# Row 26: InternetComputer(DFINITY)
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.91'
$ws.Range('E26').Value = '  +0.66%  '

# Row 27: PEPE
$ws.Range('D27').Value = '0.0₃0976'
$ws.Range('E27').Value = '  +7.81%  '

# Row 28: Binance-PegBSC-USD
$ws.Range('E28').Value = '  +0.49%  '

# Row 29: RenderToken
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.68'
$ws.Range('E29').Value = '  +7.88%  '

# Row 30: PancakeSwap
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.00'
$ws.Range('E30').Value = '  +4.60%  '

# Row 31: NEARProtocol
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.67'
$ws.Range('E31').Value = '  +5.03%  '

# Row 32: EthereumClassic
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '23.13'
$ws.Range('E32').Value = '  +3.13%  '

# Row 33: Fetch.AI
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.32'
$ws.Range('E33').Value = '  +7.07%  '

# Row 34: Aptos
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.74'
$ws.Range('E34').Value = '  +1.41%  '

# Row 35: Monero -> ImmutableX
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  +9.27%  '

# Row 36: ImmutableX -> Monero
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '160.63'
$ws.Range('E36').Value = '  +2.66%  '

# Row 37: Stacks
$ws.Range('E37').Value = '  +12.40%  '

# Row 38: EnergySwap
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '27.13'
$ws.Range('E38').Value = '  +5.25%  '

# Row 39: Hedera
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0745'
$ws.Range('E39').Value = '  +4.52%  '

# Row 40: Maker
$ws.Range('D40').Value = '2.837.21'
$ws.Range('E40').Value = '  +1.76%  '

# Row 41: VeChain
$ws.Range('E41').Value = '  +8.51%  '

# Row 42: Filecoin
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.35'
$ws.Range('E42').Value = '  +2.28%  '

# Row 43: Mantle
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.751'
$ws.Range('E43').Value = '  +2.86%  '

# Row 44: OKB
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '40.79'
$ws.Range('E44').Value = '  +2.51%  '

# Row 45: ONDO
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.06'
$ws.Range('E45').Value = '  +5.57%  '

# Row 46: InjectiveProtocol
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '22.25'
$ws.Range('E46').Value = '  +8.21%  '

# Row 47: RenzoRestakedETH
$ws.Range('D47').Value = '3.403.59'
$ws.Range('E47').Value = '  +5.20%  '

# Row 48: Stellar
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.104'
$ws.Range('E48').Value = '  +2.42%  '

# Row 49: Cosmos
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.33'
$ws.Range('E49').Value = '  +2.39%  '

# Row 50: SuiNetwork
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.812'
$ws.Range('E50').Value = '  +0.93%  '

# Row 51: Bittensor
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '283.43'
$ws.Range('E51').Value = '  +6.92%  '

